$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New "Comment" column (F) -------------------------------------------------
# ---------------------------------------------------------------------------
# Column width (engine quantizes ColumnWidth input to ~1/6 char steps, so this
# is the closest achievable approximation of the authored 47.85546875 width).
$ws.Columns.Item(6).ColumnWidth = 47

# Header cell - same bold/centered look as the rest of row 1, reuse E1's style.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "Comment"

# ---------------------------------------------------------------------------
# 2. Fix the typo in the "copy paste" review point (row 7, column D) --------
# ---------------------------------------------------------------------------
$d7 = $ws.Range("D7").Value2
$d7 = $d7 -replace "from CYRS please provie", "from CRS please provide"
$ws.Range("D7").Value = $d7

# ---------------------------------------------------------------------------
# 3. Review point 3 ("Project name shall be in the middle...") gets a new
#    reviewer comment and is reopened. ---------------------------------------
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "Open"

# F3 uses the same wrapped/left/top style as the other long comment cells
# (D5 already carries that exact style) - copy it across so no new style is
# minted for this one.
$ws.Range("D5").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Value = "Mali 25/1/2020: `nPlease remove `"1.Project Name`" no need for it"

# ---------------------------------------------------------------------------
# 4. Close out the remaining reviewed points and attach the standard
#    "reviewed and closed" comment. -------------------------------------------
# ---------------------------------------------------------------------------
$closedComment = "Mali 25/1/2020: Point is reviewed and closed"

# Row 4 - create the new left/top (no wrap) style once here ...
$ws.Range("E4").Value = "Closed"
$ws.Range("F4").Value = $closedComment
$ws.Range("F4").VerticalAlignment = -4160
$ws.Range("F4").HorizontalAlignment = -4131

# ... then reuse that exact style for rows 6 and 11 so the style table stays
# minimal (matches the single new cellXfs entry in the target workbook).
$ws.Range("F4").Copy() | Out-Null

$ws.Range("E6").Value = "Closed"
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("F6").Value = $closedComment

$ws.Range("E11").Value = "Closed"
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Value = $closedComment

# ---------------------------------------------------------------------------
# 5. Selection / scroll position left by the author after the edit ----------
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("E9").Select() | Out-Null

Write-Output "edit applied"
